$wb = $excel.ActiveWorkbook

# Insert a new worksheet for "2022-Q4" positioned after "总计" and before "2022-Q1"
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$newSheet = $wb.Worksheets.Add($wsQ1)
$newSheet.Name = "2022-Q4"

# Populate header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"
$newSheet.Range("B1:H1").Style = $wsQ1.Range("B1:H1").Style

# Row 2 data
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "011205"
$newSheet.Cells.Item(2, 3).Value = "兴银中证500指数增强C"
$newSheet.Cells.Item(2, 4).Value = "0.66"
$newSheet.Cells.Item(2, 5).Value = "84.84"
$newSheet.Cells.Item(2, 6).Value = "0.97"
$newSheet.Cells.Item(2, 7).Value = "0.0064"
$newSheet.Cells.Item(2, 8).Value = 2
$newSheet.Range("A2").Style = $wsQ1.Range("A2").Style

# Row 3 data
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "010253"
$newSheet.Cells.Item(3, 3).Value = "兴银中证500指数增强A"
$newSheet.Cells.Item(3, 4).Value = "0.32"
$newSheet.Cells.Item(3, 5).Value = "84.84"
$newSheet.Cells.Item(3, 6).Value = "0.97"
$newSheet.Cells.Item(3, 7).Value = "0.0031"
$newSheet.Cells.Item(3, 8).Value = 2
$newSheet.Range("A3").Style = $wsQ1.Range("A3").Style

# Now update the "总计" sheet: insert a new row 2 for 2022-Q4, shifting existing rows down
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.01
$wsTotal.Range("A2").Style = $wsTotal.Range("A3").Style

# Update A column index values for rows 3 and 4 (was row2=>A2=0, row3=>A3=1 before insert; now row3 and row4)
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2
